# Add new FEBRERO (February) data rows (142-168) to Sheet1
# and update the active selection to match the author's final cursor position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 142
$ws.Cells.Item(142, 1).Value = 2026
$ws.Cells.Item(142, 2).Value = "FEBRERO"
$ws.Cells.Item(142, 3).Value = 6
$ws.Cells.Item(142, 4).Value = "AMARILLO"
$ws.Cells.Item(142, 5).Value = "COLORES"
$ws.Cells.Item(142, 6).Value = "GOLDFINCH"
$ws.Cells.Item(142, 7).Value = 7295

# Row 143
$ws.Cells.Item(143, 1).Value = 2026
$ws.Cells.Item(143, 2).Value = "FEBRERO"
$ws.Cells.Item(143, 3).Value = 6
$ws.Cells.Item(143, 4).Value = "AMARILLO"
$ws.Cells.Item(143, 5).Value = "COLORES"
$ws.Cells.Item(143, 6).Value = "HIGH AND EXOTIC"
$ws.Cells.Item(143, 7).Value = 17636

# Row 144
$ws.Cells.Item(144, 1).Value = 2026
$ws.Cells.Item(144, 2).Value = "FEBRERO"
$ws.Cells.Item(144, 3).Value = 6
$ws.Cells.Item(144, 4).Value = "AMARILLO"
$ws.Cells.Item(144, 5).Value = "COLORES"
$ws.Cells.Item(144, 6).Value = "MOMENTUM"
$ws.Cells.Item(144, 7).Value = 10565

# Row 145
$ws.Cells.Item(145, 1).Value = 2026
$ws.Cells.Item(145, 2).Value = "FEBRERO"
$ws.Cells.Item(145, 3).Value = 6
$ws.Cells.Item(145, 4).Value = "AMARILLO"
$ws.Cells.Item(145, 5).Value = "COLORES"
$ws.Cells.Item(145, 6).Value = "SUNDAY MORNING"
$ws.Cells.Item(145, 7).Value = 425

# Row 146
$ws.Cells.Item(146, 1).Value = 2026
$ws.Cells.Item(146, 2).Value = "FEBRERO"
$ws.Cells.Item(146, 3).Value = 6
$ws.Cells.Item(146, 4).Value = "BIC. AMARILLO"
$ws.Cells.Item(146, 5).Value = "COLORES"
$ws.Cells.Item(146, 6).Value = "SUMMER LIGHT"
$ws.Cells.Item(146, 7).Value = 2450

# Row 147
$ws.Cells.Item(147, 1).Value = 2026
$ws.Cells.Item(147, 2).Value = "FEBRERO"
$ws.Cells.Item(147, 3).Value = 6
$ws.Cells.Item(147, 4).Value = "BICOLOR"
$ws.Cells.Item(147, 5).Value = "COLORES"
$ws.Cells.Item(147, 6).Value = "BLUSH"
$ws.Cells.Item(147, 7).Value = 2103

# Row 148
$ws.Cells.Item(148, 1).Value = 2026
$ws.Cells.Item(148, 2).Value = "FEBRERO"
$ws.Cells.Item(148, 3).Value = 6
$ws.Cells.Item(148, 4).Value = "BICOLOR"
$ws.Cells.Item(148, 5).Value = "COLORES"
$ws.Cells.Item(148, 6).Value = "DISCOVERY"
$ws.Cells.Item(148, 7).Value = 1730

# Row 149
$ws.Cells.Item(149, 1).Value = 2026
$ws.Cells.Item(149, 2).Value = "FEBRERO"
$ws.Cells.Item(149, 3).Value = 6
$ws.Cells.Item(149, 4).Value = "BLANCO"
$ws.Cells.Item(149, 5).Value = "COLORES"
$ws.Cells.Item(149, 6).Value = "HIGH AND PURE"
$ws.Cells.Item(149, 7).Value = 19790

# Row 150
$ws.Cells.Item(150, 1).Value = 2026
$ws.Cells.Item(150, 2).Value = "FEBRERO"
$ws.Cells.Item(150, 3).Value = 6
$ws.Cells.Item(150, 4).Value = "BLANCO"
$ws.Cells.Item(150, 5).Value = "COLORES"
$ws.Cells.Item(150, 6).Value = "SUGGAR DOLL"
$ws.Cells.Item(150, 7).Value = 4550

# Row 151
$ws.Cells.Item(151, 1).Value = 2026
$ws.Cells.Item(151, 2).Value = "FEBRERO"
$ws.Cells.Item(151, 3).Value = 6
$ws.Cells.Item(151, 4).Value = "BLANCO"
$ws.Cells.Item(151, 5).Value = "COLORES"
$ws.Cells.Item(151, 6).Value = "VANILLA ICE"
$ws.Cells.Item(151, 7).Value = 9078

# Row 152
$ws.Cells.Item(152, 1).Value = 2026
$ws.Cells.Item(152, 2).Value = "FEBRERO"
$ws.Cells.Item(152, 3).Value = 6
$ws.Cells.Item(152, 4).Value = "BLANCO"
$ws.Cells.Item(152, 5).Value = "COLORES"
$ws.Cells.Item(152, 6).Value = "VENDELA"
$ws.Cells.Item(152, 7).Value = 7055

# Row 153
$ws.Cells.Item(153, 1).Value = 2026
$ws.Cells.Item(153, 2).Value = "FEBRERO"
$ws.Cells.Item(153, 3).Value = 6
$ws.Cells.Item(153, 4).Value = "DURAZNO"
$ws.Cells.Item(153, 5).Value = "COLORES"
$ws.Cells.Item(153, 6).Value = "TIFANY"
$ws.Cells.Item(153, 7).Value = 3556

# Row 154
$ws.Cells.Item(154, 1).Value = 2026
$ws.Cells.Item(154, 2).Value = "FEBRERO"
$ws.Cells.Item(154, 3).Value = 6
$ws.Cells.Item(154, 4).Value = "HOT PINK"
$ws.Cells.Item(154, 5).Value = "COLORES"
$ws.Cells.Item(154, 6).Value = "COTTON CANDY"
$ws.Cells.Item(154, 7).Value = 5880

# Row 155
$ws.Cells.Item(155, 1).Value = 2026
$ws.Cells.Item(155, 2).Value = "FEBRERO"
$ws.Cells.Item(155, 3).Value = 6
$ws.Cells.Item(155, 4).Value = "HOT PINK"
$ws.Cells.Item(155, 5).Value = "COLORES"
$ws.Cells.Item(155, 6).Value = "JACARANDA"
$ws.Cells.Item(155, 7).Value = 18499

# Row 156
$ws.Cells.Item(156, 1).Value = 2026
$ws.Cells.Item(156, 2).Value = "FEBRERO"
$ws.Cells.Item(156, 3).Value = 6
$ws.Cells.Item(156, 4).Value = "HOT PINK"
$ws.Cells.Item(156, 5).Value = "COLORES"
$ws.Cells.Item(156, 6).Value = "PINK FLOYD"
$ws.Cells.Item(156, 7).Value = 13280

# Row 157
$ws.Cells.Item(157, 1).Value = 2026
$ws.Cells.Item(157, 2).Value = "FEBRERO"
$ws.Cells.Item(157, 3).Value = 6
$ws.Cells.Item(157, 4).Value = "LAVANDER"
$ws.Cells.Item(157, 5).Value = "COLORES"
$ws.Cells.Item(157, 6).Value = "DEEP PURPLE"
$ws.Cells.Item(157, 7).Value = 26430

# Row 158
$ws.Cells.Item(158, 1).Value = 2026
$ws.Cells.Item(158, 2).Value = "FEBRERO"
$ws.Cells.Item(158, 3).Value = 6
$ws.Cells.Item(158, 4).Value = "LAVANDER"
$ws.Cells.Item(158, 5).Value = "COLORES"
$ws.Cells.Item(158, 6).Value = "MOODY BLUES"
$ws.Cells.Item(158, 7).Value = 12143

# Row 159
$ws.Cells.Item(159, 1).Value = 2026
$ws.Cells.Item(159, 2).Value = "FEBRERO"
$ws.Cells.Item(159, 3).Value = 6
$ws.Cells.Item(159, 4).Value = "NARANJA"
$ws.Cells.Item(159, 5).Value = "COLORES"
$ws.Cells.Item(159, 6).Value = "ALIVE"
$ws.Cells.Item(159, 7).Value = 2675

# Row 160
$ws.Cells.Item(160, 1).Value = 2026
$ws.Cells.Item(160, 2).Value = "FEBRERO"
$ws.Cells.Item(160, 3).Value = 6
$ws.Cells.Item(160, 4).Value = "NARANJA"
$ws.Cells.Item(160, 5).Value = "COLORES"
$ws.Cells.Item(160, 6).Value = "BROMO"
$ws.Cells.Item(160, 7).Value = 4542

# Row 161
$ws.Cells.Item(161, 1).Value = 2026
$ws.Cells.Item(161, 2).Value = "FEBRERO"
$ws.Cells.Item(161, 3).Value = 6
$ws.Cells.Item(161, 4).Value = "NARANJA"
$ws.Cells.Item(161, 5).Value = "COLORES"
$ws.Cells.Item(161, 6).Value = "CLEMENTINA"
$ws.Cells.Item(161, 7).Value = 4360

# Row 162
$ws.Cells.Item(162, 1).Value = 2026
$ws.Cells.Item(162, 2).Value = "FEBRERO"
$ws.Cells.Item(162, 3).Value = 6
$ws.Cells.Item(162, 4).Value = "NARANJA"
$ws.Cells.Item(162, 5).Value = "COLORES"
$ws.Cells.Item(162, 6).Value = "NINA"
$ws.Cells.Item(162, 7).Value = 29301

# Row 163
$ws.Cells.Item(163, 1).Value = 2026
$ws.Cells.Item(163, 2).Value = "FEBRERO"
$ws.Cells.Item(163, 3).Value = 6
$ws.Cells.Item(163, 4).Value = "ROJO"
$ws.Cells.Item(163, 5).Value = "ROJO"
$ws.Cells.Item(163, 6).Value = "FREEDOM"
$ws.Cells.Item(163, 7).Value = 177428

# Row 164
$ws.Cells.Item(164, 1).Value = 2026
$ws.Cells.Item(164, 2).Value = "FEBRERO"
$ws.Cells.Item(164, 3).Value = 6
$ws.Cells.Item(164, 4).Value = "ROSADO"
$ws.Cells.Item(164, 5).Value = "COLORES"
$ws.Cells.Item(164, 6).Value = "ABSOLUT IN PINK"
$ws.Cells.Item(164, 7).Value = 3455

# Row 165
$ws.Cells.Item(165, 1).Value = 2026
$ws.Cells.Item(165, 2).Value = "FEBRERO"
$ws.Cells.Item(165, 3).Value = 6
$ws.Cells.Item(165, 4).Value = "ROSADO"
$ws.Cells.Item(165, 5).Value = "COLORES"
$ws.Cells.Item(165, 6).Value = "HIGH AND BONITA"
$ws.Cells.Item(165, 7).Value = 11032

# Row 166
$ws.Cells.Item(166, 1).Value = 2026
$ws.Cells.Item(166, 2).Value = "FEBRERO"
$ws.Cells.Item(166, 3).Value = 6
$ws.Cells.Item(166, 4).Value = "ROSADO"
$ws.Cells.Item(166, 5).Value = "COLORES"
$ws.Cells.Item(166, 6).Value = "LUCIANO"
$ws.Cells.Item(166, 7).Value = 1646

# Row 167
$ws.Cells.Item(167, 1).Value = 2026
$ws.Cells.Item(167, 2).Value = "FEBRERO"
$ws.Cells.Item(167, 3).Value = 6
$ws.Cells.Item(167, 4).Value = "ROSADO"
$ws.Cells.Item(167, 5).Value = "COLORES"
$ws.Cells.Item(167, 6).Value = "STARFISH"
$ws.Cells.Item(167, 7).Value = 5168

# Row 168
$ws.Cells.Item(168, 1).Value = 2026
$ws.Cells.Item(168, 2).Value = "FEBRERO"
$ws.Cells.Item(168, 3).Value = 6
$ws.Cells.Item(168, 4).Value = "ROSADO"
$ws.Cells.Item(168, 5).Value = "COLORES"
$ws.Cells.Item(168, 6).Value = "TABATHA"
$ws.Cells.Item(168, 7).Value = 11415

# Match the final cursor/selection state recorded in the workbook view
$null = $ws.Range("A168").Select()
